$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44435
$ws.Cells.Item(2, 10).Value = 140
$ws.Cells.Item(2, 11).Value = 21000
$ws.Cells.Item(2, 12).Value = 23000
$ws.Cells.Item(2, 13).Value = 21714
$ws.Cells.Item(2, 16).Value = 1448

# Row 3
$ws.Cells.Item(3, 4).Value = 44391
$ws.Cells.Item(3, 10).Value = 160
$ws.Cells.Item(3, 11).Value = 20000
$ws.Cells.Item(3, 12).Value = 20000
$ws.Cells.Item(3, 13).Value = 20000
$ws.Cells.Item(3, 16).Value = 1333

# Row 4
$ws.Cells.Item(4, 4).Value = 44392
$ws.Cells.Item(4, 10).Value = 220
$ws.Cells.Item(4, 11).Value = 23000
$ws.Cells.Item(4, 12).Value = 23000
$ws.Cells.Item(4, 13).Value = 23000
$ws.Cells.Item(4, 16).Value = 1533

# Row 5
$ws.Cells.Item(5, 4).Value = 44406
$ws.Cells.Item(5, 10).Value = 400
$ws.Cells.Item(5, 11).Value = 20000
$ws.Cells.Item(5, 12).Value = 22000
$ws.Cells.Item(5, 13).Value = 20850
$ws.Cells.Item(5, 16).Value = 1390

# Row 6
$ws.Cells.Item(6, 4).Value = 44476
$ws.Cells.Item(6, 10).Value = 220
$ws.Cells.Item(6, 11).Value = 20000
$ws.Cells.Item(6, 12).Value = 22000
$ws.Cells.Item(6, 13).Value = 20909
$ws.Cells.Item(6, 16).Value = 1394

# Row 7
$ws.Cells.Item(7, 4).Value = 44446
$ws.Cells.Item(7, 10).Value = 150
$ws.Cells.Item(7, 11).Value = 22000
$ws.Cells.Item(7, 12).Value = 24000
$ws.Cells.Item(7, 13).Value = 22667
$ws.Cells.Item(7, 16).Value = 1511

# Row 8
$ws.Cells.Item(8, 4).Value = 44398
$ws.Cells.Item(8, 10).Value = 130
$ws.Cells.Item(8, 11).Value = 20000
$ws.Cells.Item(8, 12).Value = 20000
$ws.Cells.Item(8, 13).Value = 20000
$ws.Cells.Item(8, 16).Value = 1333

# Row 9
$ws.Cells.Item(9, 4).Value = 44449
$ws.Cells.Item(9, 10).Value = 220
$ws.Cells.Item(9, 11).Value = 22000
$ws.Cells.Item(9, 12).Value = 24000
$ws.Cells.Item(9, 13).Value = 23091
$ws.Cells.Item(9, 16).Value = 1539

# Row 10
$ws.Cells.Item(10, 4).Value = 44483
$ws.Cells.Item(10, 10).Value = 220
$ws.Cells.Item(10, 11).Value = 18000
$ws.Cells.Item(10, 12).Value = 20000
$ws.Cells.Item(10, 13).Value = 18909
$ws.Cells.Item(10, 16).Value = 1261

# Row 12
$ws.Cells.Item(12, 4).Value = 44453
$ws.Cells.Item(12, 10).Value = 280
$ws.Cells.Item(12, 11).Value = 20000
$ws.Cells.Item(12, 12).Value = 22000
$ws.Cells.Item(12, 13).Value = 21286
$ws.Cells.Item(12, 16).Value = 1419

# Row 14
$ws.Cells.Item(14, 4).Value = 44400
$ws.Cells.Item(14, 10).Value = 130
$ws.Cells.Item(14, 11).Value = 24000
$ws.Cells.Item(14, 12).Value = 24000
$ws.Cells.Item(14, 13).Value = 24000
$ws.Cells.Item(14, 16).Value = 1600
